$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Add two new list paragraphs after "Used for playing the game"
#    ("Random Encounters!" and "Travel! "), matching the existing
#    level-1 / numId-2 list formatting used throughout that section.
# ------------------------------------------------------------------
$target = $null
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Used for playing the game*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
    $p1 = $target.Next()
    $p1.Range.Text = "Random Encounters!"

    $p1.Range.InsertParagraphAfter()
    $p2 = $p1.Next()
    $p2.Range.Text = "Travel! "
}

# ------------------------------------------------------------------
# 2. Register the new ListLabel82..ListLabel99 character styles that
#    back the additional list levels introduced by the restructure.
# ------------------------------------------------------------------
$listLabelFonts = @(
    @{ Id = "ListLabel82"; Font = "Symbol" },
    @{ Id = "ListLabel83"; Font = "Courier New" },
    @{ Id = "ListLabel84"; Font = "Wingdings" },
    @{ Id = "ListLabel85"; Font = "Symbol" },
    @{ Id = "ListLabel86"; Font = "Courier New" },
    @{ Id = "ListLabel87"; Font = "Wingdings" },
    @{ Id = "ListLabel88"; Font = "Symbol" },
    @{ Id = "ListLabel89"; Font = "Courier New" },
    @{ Id = "ListLabel90"; Font = "Wingdings" },
    @{ Id = "ListLabel91"; Font = "OpenSymbol" },
    @{ Id = "ListLabel92"; Font = "OpenSymbol" },
    @{ Id = "ListLabel93"; Font = "OpenSymbol" },
    @{ Id = "ListLabel94"; Font = "OpenSymbol" },
    @{ Id = "ListLabel95"; Font = "OpenSymbol" },
    @{ Id = "ListLabel96"; Font = "OpenSymbol" },
    @{ Id = "ListLabel97"; Font = "OpenSymbol" },
    @{ Id = "ListLabel98"; Font = "OpenSymbol" },
    @{ Id = "ListLabel99"; Font = "OpenSymbol" }
)

foreach ($entry in $listLabelFonts) {
    $num = $entry.Id.Substring(9)
    $style = $d.Styles.Add($entry.Id, 2)
    $style.NameLocal = "ListLabel " + $num
    $style.QuickStyle = $true
    $style.Font.NameBi = $entry.Font
}
